$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E2").Value = 11.14

$ws.Range("E3").Value = 10.42

$ws.Range("E4").Value = 10.02
$ws.Range("F4").Value = 10.31

$ws.Range("B5").Value = 8.859999999999999
$ws.Range("C5").Value = 9.58
$ws.Range("D5").Value = 9.98
$ws.Range("F5").Value = 10.27
$ws.Range("G5").Value = 9.23
$ws.Range("H5").Value = 7.74

$ws.Range("D6").Value = 9.69
$ws.Range("E6").Value = 9.73
$ws.Range("G6").Value = 10.39
$ws.Range("H6").Value = 11.37
$ws.Range("J6").Value = 8.18

$ws.Range("E7").Value = 10.77
$ws.Range("F7").Value = 9.609999999999999

$ws.Range("E8").Value = 12.26
$ws.Range("F8").Value = 8.630000000000001
$ws.Range("I8").Value = 7.57

$ws.Range("H9").Value = 12.43

$ws.Range("F10").Value = 11.82
